$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new match rows (Chinese Super League) were inserted after row 2,
# pushing the former rows 3-8 down to rows 5-10. Insert the rows first so
# row numbers below line up with the final layout (dimension becomes A1:AO10).
$ws.Rows.Item(3).Resize(2).Insert()

# Columns B (Date) and C (Time) hold plain text like "2025-10-31" / "05:35:00".
# Force text format first so Excel does not auto-convert them into date/time serials.
$ws.Range("B2:C10").NumberFormat = "@"

$rowData = @{}
$rowData[2] = @("Australian A-League Men", "2025-10-31", "05:35:00", "Brisbane Roar", "Melbourne City", 4, 4.4, 1.94, 2.04, 3.9, 4, 1.43, 1.07, 3.5, 1.35, 1.85, 1.99, 1.34, 3.6, 1.87, 1.96, 1.97, 1.29, 16, 9, 14, 1000, 15, 8.6, 11, 24, 1000, 21, 20, 1000, 120, 70, 1000, 140, 1000, 16.5)
$rowData[3] = @("Chinese Super League", "2025-10-31", "07:00:00", "Shanghai Port FC", "Zhejiang Greentown", 1.53, 1.63, 4.3, 6.4, 4.4, 6.4, 1.14, 1.01, 1.01, 1.07, 4, 1.2, 2.36, 1.5, 1.36, 3.2, 1.2, 2.56, 75, 55, 65, 130, 980, 980, 980, 60, 980, 980, 21, 980, 980, 980, 980, 980, 3.9, 980)
$rowData[4] = @("Australian A-League Men", "2025-10-31", "07:45:00", "Perth Glory", "Melbourne Victory", 3.75, 3.95, 1.92, 2.02, 4.1, 4.5, 1.29, 1.03, 5.4, 1.2, 2.5, 1.59, 1.61, 2.46, 1.58, 2.5, 1.99, 1.33, 30, 16.5, 18.5, 28, 21, 12.5, 13.5, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 11)
$rowData[5] = @("Chinese Super League", "2025-10-31", "08:35:00", "Chengdu Rongcheng", "Henan Songshan Longmen", 1.35, 1.49, 8.2, 9.2, 4.9, 6.8, 0, 0, 0, 0, 3.45, 1.32, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[6] = @("Saudi Professional League", "2025-10-31", "10:40:00", "Al-Akhdoud", "Al Najma Club", 2.2, 2.68, 3.15, 4.1, 3.25, 3.95, 0, 0, 0, 0, 1.87, 1.81, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[7] = @("Saudi Professional League", "2025-10-31", "11:50:00", "Al-Hilal", "Al-Shabab (KSA)", 1.24, 1.3, 11.5, 15.5, 6.6, 8, 0, 0, 0, 0, 2.86, 1.42, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[8] = @("Saudi Professional League", "2025-10-31", "14:30:00", "Al-Ettifaq", "Al-Hazm (KSA)", 2, 2.16, 3.75, 4.3, 3.5, 4.1, 0, 0, 0, 0, 1.97, 1.83, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[9] = @("German Bundesliga", "2025-10-31", "16:30:00", "Augsburg", "Dortmund", 4.9, 5.2, 1.7, 1.74, 4.3, 4.6, 0, 1.03, 5.9, 1.18, 2.66, 1.54, 1.69, 2.32, 1.6, 2.52, 0, 0, 28, 14, 13, 19.5, 28, 10.5, 10.5, 15.5, 48, 21, 17, 27, 130, 65, 980, 65, 42, 6.8)
$rowData[10] = @("Spanish La Liga", "2025-10-31", "17:00:00", "Getafe", "Girona", 2.26, 2.28, 4.1, 4.4, 3.05, 3.1, 0, 1.14, 2.68, 1.58, 1.54, 2.74, 1.19, 5.8, 2.24, 1.76, 0, 0, 8.2, 11, 29, 120, 7, 7.2, 18.5, 80, 12.5, 12.5, 26, 110, 32, 34, 70, 1000, 32, 150)

foreach ($r in $rowData.Keys) {
    $values = $rowData[$r]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $values[$c]
    }
}

Write-Host "Applied edits"
